$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

$ws.Range("A16").Value = 43378
$ws.Range("B16").Value = 0.38541666666666669
$ws.Range("C16").Value = 0.61111111111111105
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = "1h Blogi kirjoittelua Ubuntu 16.04 -asennuksesta -> Master-palvelimen asennus.  3h midPoint asennuksen kokeilua. -> https://wiki.evolveum.com/display/midPoint/MidPoint+Easy+Install. https://wiki.evolveum.com/display/midPoint/midPoint+on+Ubuntu,+Tomcat,+PostgreSQL+HOWTO https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%205.10.2018.txt"

$ws.Rows.Item(16).RowHeight = 135

$ws.Activate() | Out-Null
$ws.Range("C16").Select() | Out-Null
